$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove existing hyperlinks (Range.Hyperlinks.Delete clears all sheet hyperlinks) ---
$ws.Range("A1").Hyperlinks.Delete()

# --- Clear old data rows (rows 2-3) before rewriting the full data set ---
$ws.Range("A2:H3").ClearContents()

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 50 - 5/6
$ws.Columns.Item(4).ColumnWidth = 30 - 5/6
$ws.Columns.Item(8).ColumnWidth = 18 - 5/6

# --- Row 2: new PM job posting ---
$ws.Range("A2").Value = "2025-10-05 12:32:23"
$ws.Range("B2").Value = "【PM募集】生成AIを活用した新規業務改善SaaS開発プロジェクトのプロジェクトマネージャー"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5407076"
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = "🔥AI,Ai ◆開発 ◇業務改善"

# --- Row 3: new WEB system job posting ---
$ws.Range("A3").Value = "2025-10-05 12:32:23"
$ws.Range("B3").Value = "画像に情報を紐づけるWEBシステムの開発"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5407101"
$ws.Range("G3").Value = 78
$ws.Range("H3").Value = "◆開発"

# --- Row 4: new Power Automate job posting ---
$ws.Range("A4").Value = "2025-10-05 12:32:23"
$ws.Range("B4").Value = "【Power Automate for Desktop】販売管理システムへExcelから自動入力"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5407216"
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = "◇管理"

# --- Row 5: new FX EA job posting (different price/url than row 6) ---
$ws.Range("A5").Value = "2025-10-05 12:32:23"
$ws.Range("B5").Value = "【急募】FXトレード履歴を基にしたEA作成依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5407099"
$ws.Range("G5").Value = 18

# --- Row 6: previously row 2, timestamp refreshed ---
$ws.Range("A6").Value = "2025-10-05 12:32:23"
$ws.Range("B6").Value = "【急募】FXトレード履歴を基にしたEA作成依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5406904"
$ws.Range("G6").Value = 18

# --- Row 7: previously row 3, timestamp refreshed ---
$ws.Range("A7").Value = "2025-10-05 12:32:23"
$ws.Range("B7").Value = "【SalesIQ活用】CRMと連携したリード獲得方法を教えてください"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "~ 5,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5400402"
$ws.Range("G7").Value = 10

# --- Hyperlinks for URL column, in row order ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5407076")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5407101")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5407216")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5407099")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5406904")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5400402")
